# Updates the "cryptos" price table: new Price (column D) and Volume(1h)
# (column E) figures for the latest GitHub Actions refresh.
#
# Column D cells hold numeric-looking text (e.g. "232.46", "19.00",
# "0.0617") that must stay literal text exactly as scraped. Assigning such
# strings straight to .Value lets Excel auto-coerce them to real numbers
# (dropping trailing zeros, flipping to scientific notation, etc.), so each
# D cell is forced to text via NumberFormat "@" before the write, then
# restored to the default "Normal" style so it doesn't pick up a stray
# explicit style index.
#
# Column E values are percentage strings with deliberate leading/trailing
# padding spaces (e.g. "  +0.25%  ") that must be preserved verbatim;
# assigning plain strings there is safe since they never look numeric.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.269.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.059.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("E6").Value = "  +2.51%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.79"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.362.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.776"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.060.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.260.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0808"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "226.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.36%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.17%  "
$ws.Range("E28").Value = "  +6.81%  "
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("E30").Value = "  -2.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0617"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.66%  "
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("E39").Value = "  -2.32%  "
$ws.Range("E40").Value = "  -4.21%  "
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.472.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0942"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0212"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.37%  "
$ws.Range("E46").Value = "  +3.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "15.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("E51").Value = "  +0.74%  "
